$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header style (bold, border, center) from A1 to the new header columns G1:AE1
$ws.Range("A1").Copy($ws.Range("G1:AE1"))

# Row 1
$ws.Cells.Item(1, 7).Value = "ARAP-InRays Initial VS Mov (%)"
$ws.Cells.Item(1, 8).Value = "ARAP-TwoPoints Improvement (%)"
$ws.Cells.Item(1, 9).Value = "ARAP-TwoPoints Final Vs Mov (%)"
$ws.Cells.Item(1, 10).Value = "ARAP-TwoPoints Initial VS Mov (%)"
$ws.Cells.Item(1, 11).Value = "ARAP-FarPoints Improvement (%)"
$ws.Cells.Item(1, 12).Value = "ARAP-FarPoints Final Vs Mov (%)"
$ws.Cells.Item(1, 13).Value = "ARAP-FarPoints Initial VS Mov (%)"
$ws.Cells.Item(1, 14).Value = "ARAP_depth-InRays Improvement (%)"
$ws.Cells.Item(1, 15).Value = "ARAP_depth-InRays Final Vs Mov (%)"
$ws.Cells.Item(1, 16).Value = "ARAP_depth-InRays Initial VS Mov (%)"
$ws.Cells.Item(1, 17).Value = "ARAP_depth-TwoPoints Improvement (%)"
$ws.Cells.Item(1, 18).Value = "ARAP_depth-TwoPoints Final Vs Mov (%)"
$ws.Cells.Item(1, 19).Value = "ARAP_depth-TwoPoints Initial VS Mov (%)"
$ws.Cells.Item(1, 20).Value = "ARAP_depth-FarPoints Improvement (%)"
$ws.Cells.Item(1, 21).Value = "ARAP_depth-FarPoints Final Vs Mov (%)"
$ws.Cells.Item(1, 22).Value = "ARAP_depth-FarPoints Initial VS Mov (%)"
$ws.Cells.Item(1, 23).Value = "ARAP_depth_onlyTriang-InRays Improvement (%)"
$ws.Cells.Item(1, 24).Value = "ARAP_depth_onlyTriang-InRays Final Vs Mov (%)"
$ws.Cells.Item(1, 25).Value = "ARAP_depth_onlyTriang-InRays Initial VS Mov (%)"
$ws.Cells.Item(1, 26).Value = "ARAP_depth_onlyTriang-TwoPoints Improvement (%)"
$ws.Cells.Item(1, 27).Value = "ARAP_depth_onlyTriang-TwoPoints Final Vs Mov (%)"
$ws.Cells.Item(1, 28).Value = "ARAP_depth_onlyTriang-TwoPoints Initial VS Mov (%)"
$ws.Cells.Item(1, 29).Value = "ARAP_depth_onlyTriang-FarPoints Improvement (%)"
$ws.Cells.Item(1, 30).Value = "ARAP_depth_onlyTriang-FarPoints Final Vs Mov (%)"
$ws.Cells.Item(1, 31).Value = "ARAP_depth_onlyTriang-FarPoints Initial VS Mov (%)"

# Row 2
$ws.Cells.Item(2, 7).Value = 77.97
$ws.Cells.Item(2, 8).Value = 4.25
$ws.Cells.Item(2, 9).Value = 81.78
$ws.Cells.Item(2, 10).Value = 85.31999999999999
$ws.Cells.Item(2, 11).Value = -5.87
$ws.Cells.Item(2, 12).Value = 91.40000000000001
$ws.Cells.Item(2, 13).Value = 86.33
$ws.Cells.Item(2, 14).Value = -204.5
$ws.Cells.Item(2, 15).Value = 82.81999999999999
$ws.Cells.Item(2, 16).Value = 27.09
$ws.Cells.Item(2, 17).Value = -47.81
$ws.Cells.Item(2, 18).Value = 80.34999999999999
$ws.Cells.Item(2, 19).Value = 54.43
$ws.Cells.Item(2, 20).Value = -52.51
$ws.Cells.Item(2, 21).Value = 100.95
$ws.Cells.Item(2, 22).Value = 66.08
$ws.Cells.Item(2, 23).Value = -207.24
$ws.Cells.Item(2, 24).Value = 83.56
$ws.Cells.Item(2, 25).Value = 27.09
$ws.Cells.Item(2, 26).Value = -44.02
$ws.Cells.Item(2, 27).Value = 78.29000000000001
$ws.Cells.Item(2, 28).Value = 54.43
$ws.Cells.Item(2, 29).Value = -53.05
$ws.Cells.Item(2, 30).Value = 101.3
$ws.Cells.Item(2, 31).Value = 66.08

# Row 3
$ws.Cells.Item(3, 5).Value = 6.28
$ws.Cells.Item(3, 6).Value = 56.77
$ws.Cells.Item(3, 7).Value = 60.4
$ws.Cells.Item(3, 8).Value = 11.46
$ws.Cells.Item(3, 9).Value = 58.18
$ws.Cells.Item(3, 10).Value = 65.59999999999999
$ws.Cells.Item(3, 11).Value = 24.38
$ws.Cells.Item(3, 12).Value = 53.45
$ws.Cells.Item(3, 13).Value = 70.8
$ws.Cells.Item(3, 14).Value = 23.31
$ws.Cells.Item(3, 15).Value = 32.9
$ws.Cells.Item(3, 16).Value = 42.8
$ws.Cells.Item(3, 17).Value = 28.37
$ws.Cells.Item(3, 18).Value = 40.69
$ws.Cells.Item(3, 19).Value = 56.8
$ws.Cells.Item(3, 20).Value = 55.51
$ws.Cells.Item(3, 21).Value = 37.2
$ws.Cells.Item(3, 22).Value = 83.59999999999999
$ws.Cells.Item(3, 23).Value = 23.82
$ws.Cells.Item(3, 24).Value = 32.68
$ws.Cells.Item(3, 25).Value = 42.8
$ws.Cells.Item(3, 26).Value = 28.47
$ws.Cells.Item(3, 27).Value = 40.63
$ws.Cells.Item(3, 28).Value = 56.8
$ws.Cells.Item(3, 29).Value = 55.47
$ws.Cells.Item(3, 30).Value = 37.23
$ws.Cells.Item(3, 31).Value = 83.59999999999999

# Row 4
$ws.Cells.Item(4, 5).Value = -1.78
$ws.Cells.Item(4, 6).Value = 68.26000000000001
$ws.Cells.Item(4, 7).Value = 67.02
$ws.Cells.Item(4, 8).Value = 7.44
$ws.Cells.Item(4, 9).Value = 68.36
$ws.Cells.Item(4, 10).Value = 73.78
$ws.Cells.Item(4, 11).Value = 12.15
$ws.Cells.Item(4, 12).Value = 66.38
$ws.Cells.Item(4, 13).Value = 75.48
$ws.Cells.Item(4, 14).Value = -178.72
$ws.Cells.Item(4, 15).Value = 63.8
$ws.Cells.Item(4, 16).Value = 22.83
$ws.Cells.Item(4, 17).Value = -20.02
$ws.Cells.Item(4, 18).Value = 63.65
$ws.Cells.Item(4, 19).Value = 53.07
$ws.Cells.Item(4, 20).Value = -70.81999999999999
$ws.Cells.Item(4, 21).Value = 109.69
$ws.Cells.Item(4, 22).Value = 64.06
$ws.Cells.Item(4, 23).Value = -177.44
$ws.Cells.Item(4, 24).Value = 63.51
$ws.Cells.Item(4, 25).Value = 22.83
$ws.Cells.Item(4, 26).Value = -19.08
$ws.Cells.Item(4, 27).Value = 63.15
$ws.Cells.Item(4, 28).Value = 53.07
$ws.Cells.Item(4, 29).Value = -71.23999999999999
$ws.Cells.Item(4, 30).Value = 109.95
$ws.Cells.Item(4, 31).Value = 64.06

# Row 5
$ws.Cells.Item(5, 5).Value = -6.16
$ws.Cells.Item(5, 6).Value = 80.90000000000001
$ws.Cells.Item(5, 7).Value = 76.26000000000001
$ws.Cells.Item(5, 8).Value = 5.95
$ws.Cells.Item(5, 9).Value = 80.68000000000001
$ws.Cells.Item(5, 10).Value = 85.79000000000001
$ws.Cells.Item(5, 11).Value = 1.63
$ws.Cells.Item(5, 12).Value = 84.69
$ws.Cells.Item(5, 13).Value = 86.11
$ws.Cells.Item(5, 14).Value = -1648.69
$ws.Cells.Item(5, 15).Value = 116.83
$ws.Cells.Item(5, 16).Value = 6.67
$ws.Cells.Item(5, 17).Value = -121.04
$ws.Cells.Item(5, 18).Value = 111.01
$ws.Cells.Item(5, 19).Value = 50.22
$ws.Cells.Item(5, 20).Value = -1603.6
$ws.Cells.Item(5, 21).Value = 873.0700000000001
$ws.Cells.Item(5, 22).Value = 51.28
$ws.Cells.Item(5, 23).Value = -1753.71
$ws.Cells.Item(5, 24).Value = 123.85
$ws.Cells.Item(5, 25).Value = 6.67
$ws.Cells.Item(5, 26).Value = -106.49
$ws.Cells.Item(5, 27).Value = 103.7
$ws.Cells.Item(5, 28).Value = 50.22
$ws.Cells.Item(5, 29).Value = -1598.78
$ws.Cells.Item(5, 30).Value = 870.59
$ws.Cells.Item(5, 31).Value = 51.28

# Row 6
$ws.Cells.Item(6, 5).Value = 5.15
$ws.Cells.Item(6, 6).Value = 44.11
$ws.Cells.Item(6, 7).Value = 46.5
$ws.Cells.Item(6, 8).Value = 18.19
$ws.Cells.Item(6, 9).Value = 45.57
$ws.Cells.Item(6, 10).Value = 55.7
$ws.Cells.Item(6, 11).Value = 34.95
$ws.Cells.Item(6, 12).Value = 36.55
$ws.Cells.Item(6, 13).Value = 56.2
$ws.Cells.Item(6, 14).Value = 23.21
$ws.Cells.Item(6, 15).Value = 8.08
$ws.Cells.Item(6, 16).Value = 10.5
$ws.Cells.Item(6, 17).Value = 80.31
$ws.Cells.Item(6, 18).Value = 9.92
$ws.Cells.Item(6, 19).Value = 50.4
$ws.Cells.Item(6, 20).Value = 57.33
$ws.Cells.Item(6, 21).Value = 22.48
$ws.Cells.Item(6, 22).Value = 52.7
$ws.Cells.Item(6, 23).Value = 22.98
$ws.Cells.Item(6, 24).Value = 8.109999999999999
$ws.Cells.Item(6, 25).Value = 10.5
$ws.Cells.Item(6, 26).Value = 54.02
$ws.Cells.Item(6, 27).Value = 23.16
$ws.Cells.Item(6, 28).Value = 50.4
$ws.Cells.Item(6, 29).Value = 57.6
$ws.Cells.Item(6, 30).Value = 22.34
$ws.Cells.Item(6, 31).Value = 52.7

# Row 7
$ws.Cells.Item(7, 5).Value = -2.05
$ws.Cells.Item(7, 6).Value = 68.78
$ws.Cells.Item(7, 7).Value = 67.40000000000001
$ws.Cells.Item(7, 8).Value = 9.68
$ws.Cells.Item(7, 9).Value = 68.83
$ws.Cells.Item(7, 10).Value = 76.23
$ws.Cells.Item(7, 11).Value = -54.31
$ws.Cells.Item(7, 12).Value = 117.73
$ws.Cells.Item(7, 13).Value = 76.28
$ws.Cells.Item(7, 14).Value = -5558.03
$ws.Cells.Item(7, 15).Value = 311
$ws.Cells.Item(7, 16).Value = 5.49
$ws.Cells.Item(7, 17).Value = -70.43000000000001
$ws.Cells.Item(7, 18).Value = 85.48999999999999
$ws.Cells.Item(7, 19).Value = 50.16
$ws.Cells.Item(7, 20).Value = -1042.32
$ws.Cells.Item(7, 21).Value = 577.89
$ws.Cells.Item(7, 22).Value = 50.57
$ws.Cells.Item(7, 23).Value = -5574.19
$ws.Cells.Item(7, 24).Value = 311.89
$ws.Cells.Item(7, 25).Value = 5.49
$ws.Cells.Item(7, 26).Value = -68.58
$ws.Cells.Item(7, 27).Value = 84.56
$ws.Cells.Item(7, 28).Value = 50.16
$ws.Cells.Item(7, 29).Value = -1066.05
$ws.Cells.Item(7, 30).Value = 589.9
$ws.Cells.Item(7, 31).Value = 50.57

# Row 8
$ws.Cells.Item(8, 5).Value = -0.45
$ws.Cells.Item(8, 6).Value = 53.76
$ws.Cells.Item(8, 7).Value = 53.4
$ws.Cells.Item(8, 8).Value = 14.46
$ws.Cells.Item(8, 9).Value = 53.61
$ws.Cells.Item(8, 10).Value = 62.65
$ws.Cells.Item(8, 11).Value = 40.69
$ws.Cells.Item(8, 12).Value = 39.56
$ws.Cells.Item(8, 13).Value = 66.67
$ws.Cells.Item(8, 14).Value = -25.48
$ws.Cells.Item(8, 15).Value = 41.36
$ws.Cells.Item(8, 16).Value = 33.02
$ws.Cells.Item(8, 17).Value = 30.17
$ws.Cells.Item(8, 18).Value = 39.67
$ws.Cells.Item(8, 19).Value = 56.79
$ws.Cells.Item(8, 20).Value = 47.24
$ws.Cells.Item(8, 21).Value = 39.62
$ws.Cells.Item(8, 22).Value = 75
$ws.Cells.Item(8, 23).Value = -25.84
$ws.Cells.Item(8, 24).Value = 41.48
$ws.Cells.Item(8, 25).Value = 33.02
$ws.Cells.Item(8, 26).Value = 29.88
$ws.Cells.Item(8, 27).Value = 39.83
$ws.Cells.Item(8, 28).Value = 56.79
$ws.Cells.Item(8, 29).Value = 47.23
$ws.Cells.Item(8, 30).Value = 39.63
$ws.Cells.Item(8, 31).Value = 75

# Row 9
$ws.Cells.Item(9, 5).Value = -2.38
$ws.Cells.Item(9, 6).Value = 68.70999999999999
$ws.Cells.Item(9, 7).Value = 66.98
$ws.Cells.Item(9, 8).Value = 8.130000000000001
$ws.Cells.Item(9, 9).Value = 69.17
$ws.Cells.Item(9, 10).Value = 75.19
$ws.Cells.Item(9, 11).Value = 11.83
$ws.Cells.Item(9, 12).Value = 67.13
$ws.Cells.Item(9, 13).Value = 76.12
$ws.Cells.Item(9, 14).Value = -279.55
$ws.Cells.Item(9, 15).Value = 74.98
$ws.Cells.Item(9, 16).Value = 19.78
$ws.Cells.Item(9, 17).Value = -28.81
$ws.Cells.Item(9, 18).Value = 67.06999999999999
$ws.Cells.Item(9, 19).Value = 52.05
$ws.Cells.Item(9, 20).Value = -25.81
$ws.Cells.Item(9, 21).Value = 75.98
$ws.Cells.Item(9, 22).Value = 60.26
$ws.Cells.Item(9, 23).Value = -276.27
$ws.Cells.Item(9, 24).Value = 74.33
$ws.Cells.Item(9, 25).Value = 19.78
$ws.Cells.Item(9, 26).Value = -27.46
$ws.Cells.Item(9, 27).Value = 66.37
$ws.Cells.Item(9, 28).Value = 52.05
$ws.Cells.Item(9, 29).Value = -28.46
$ws.Cells.Item(9, 30).Value = 77.58
$ws.Cells.Item(9, 31).Value = 60.26

# Row 10
$ws.Cells.Item(10, 5).Value = 45.75
$ws.Cells.Item(10, 6).Value = 23.82
$ws.Cells.Item(10, 7).Value = 43.87
$ws.Cells.Item(10, 8).Value = 39.17
$ws.Cells.Item(10, 9).Value = 33.6
$ws.Cells.Item(10, 10).Value = 55.27
$ws.Cells.Item(10, 11).Value = 44.44
$ws.Cells.Item(10, 12).Value = 30.91
$ws.Cells.Item(10, 13).Value = 55.65
$ws.Cells.Item(10, 14).Value = -113.74
$ws.Cells.Item(10, 15).Value = 21.16
$ws.Cells.Item(10, 16).Value = 9.880000000000001
$ws.Cells.Item(10, 17).Value = 53.37
$ws.Cells.Item(10, 18).Value = 23.72
$ws.Cells.Item(10, 19).Value = 50.9
$ws.Cells.Item(10, 20).Value = 68.61
$ws.Cells.Item(10, 21).Value = 16.9
$ws.Cells.Item(10, 22).Value = 53.85
$ws.Cells.Item(10, 23).Value = -113.62
$ws.Cells.Item(10, 24).Value = 21.15
$ws.Cells.Item(10, 25).Value = 9.880000000000001
$ws.Cells.Item(10, 26).Value = 55.97
$ws.Cells.Item(10, 27).Value = 22.39
$ws.Cells.Item(10, 28).Value = 50.9
$ws.Cells.Item(10, 29).Value = 68.40000000000001
$ws.Cells.Item(10, 30).Value = 17.02
$ws.Cells.Item(10, 31).Value = 53.85

# Row 11
$ws.Cells.Item(11, 5).Value = -4.99
$ws.Cells.Item(11, 6).Value = 67.73
$ws.Cells.Item(11, 7).Value = 64.52
$ws.Cells.Item(11, 8).Value = 10.76
$ws.Cells.Item(11, 9).Value = 67.38
$ws.Cells.Item(11, 10).Value = 75.48
$ws.Cells.Item(11, 11).Value = -5.8
$ws.Cells.Item(11, 12).Value = 79.65000000000001
$ws.Cells.Item(11, 13).Value = 75.27
$ws.Cells.Item(11, 14).Value = -10788.04
$ws.Cells.Item(11, 15).Value = 585.6
$ws.Cells.Item(11, 16).Value = 5.4
$ws.Cells.Item(11, 17).Value = -64.11
$ws.Cells.Item(11, 18).Value = 82.36
$ws.Cells.Item(11, 19).Value = 50.18
$ws.Cells.Item(11, 20).Value = -1066.75
$ws.Cells.Item(11, 21).Value = 589.77
$ws.Cells.Item(11, 22).Value = 50.54
$ws.Cells.Item(11, 23).Value = -10821.31
$ws.Cells.Item(11, 24).Value = 587.39
$ws.Cells.Item(11, 25).Value = 5.4
$ws.Cells.Item(11, 26).Value = -68.70999999999999
$ws.Cells.Item(11, 27).Value = 84.67
$ws.Cells.Item(11, 28).Value = 50.18
$ws.Cells.Item(11, 29).Value = -1067.6
$ws.Cells.Item(11, 30).Value = 590.1900000000001
$ws.Cells.Item(11, 31).Value = 50.54
